# Update the marksheet's correct/total marks figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row (row 11): Right column B changes from 3 to 5
$ws.Range("B11").Value = 5

# "Total" row (row 12): Right column B changes from 75 to 125
$ws.Range("B12").Value = 125

# "Total" row (row 12): Max column E text changes from "75/84" to "125/140"
$ws.Range("E12").Value = "125/140"
